# Regenerate merged AHB files
# 1) Rename the "_old" / "_new" suffixed column headers to the new
#    version-specific suffixes ("_FV2404" / "_FV2410").
# 2) Turn the A1:U87 range into a native Excel Table (ListObject).
# 3) Freeze the header row (split/freeze at row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldNames = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

$newNames = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

# Columns A-J (1-10): "*_old" -> "*_FV2404"
for ($i = 0; $i -lt $oldNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldNames[$i]
}

# Column K (11): "diff" stays as-is.
$ws.Cells.Item(1, 11).Value = "diff"

# Columns L-U (12-21): "*_new" -> "*_FV2410"
for ($i = 0; $i -lt $newNames.Count; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $newNames[$i]
}

# Create the Excel Table over the full used range, using the (already
# renamed) header row as the column names.
$tableRange = $ws.Range("A1:U87")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# Freeze panes at row 2 (keep header row 1 visible while scrolling).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
